# Adds the "Labour" dataset-description section to the document:
#   - appends "(1) " to the existing "Frequency = daily" paragraph
#   - inserts a new bold/size-32 "Labour" heading
#   - inserts the Level/Description/Number-of-time-series table
#   - inserts the trailing summary paragraphs (levels/time series/horizon/frequency)

$d = $word.ActiveDocument

# --- 1) "Frequency = daily " -> "Frequency = daily (1) " ------------------
$lastPara = $d.Paragraphs.Last
$freqRange = $lastPara.Range
$freqRange.Collapse(0)
$freqRange.InsertAfter("(1) ")

# --- 2) Append the new "Labour" section after that paragraph --------------
# Collapse a range to the very end of the document body (just before the
# final section mark) and inject the new paragraphs/table as literal OOXML,
# which keeps every element (run/paragraph/table property) faithful to the
# target markup.
$bodyEnd = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$labourXml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t>Labour</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:tblPr>
        <w:tblStyle w:val="TableGrid"/>
        <w:tblW w:w="0" w:type="auto"/>
        <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="3003"/>
        <w:gridCol w:w="3003"/>
        <w:gridCol w:w="3004"/>
      </w:tblGrid>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3003" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Level</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3003" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Description</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3004" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Number of time series</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3003" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>0</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3003" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">Total </w:t>
            </w:r>
            <w:r>
              <w:t>Employees</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3004" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>1</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3003" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>1</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3003" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Main Occupation</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3004" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>8</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3003" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>2</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3003" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Employment status (Full time, Part-time)</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3004" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>8</w:t>
            </w:r>
            <w:r>
              <w:t>*</w:t>
            </w:r>
            <w:r>
              <w:t>2</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> = </w:t>
            </w:r>
            <w:r>
              <w:t>16</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3003" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>3</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3003" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Gender (Female, Male)</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3004" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>16 * 2 = 32</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">Total number of levels = </w:t>
      </w:r>
      <w:r>
        <w:t>4</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">Total number of time series = </w:t>
      </w:r>
      <w:r>
        <w:t>57</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">Horizon = </w:t>
      </w:r>
      <w:r>
        <w:t>12</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">Frequency = </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">4 (Quarterly) </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@

$bodyEnd.InsertXML($labourXml)
